# Auto-generated PowerShell Excel COM-interop script
# Applies the odds-update diff to the active worksheet (row-by-row data refresh).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AB2").Value = 12
$ws.Range("AE2").Value = 36
$ws.Range("AF2").Value = 15.5
$ws.Range("AZ2").Value = 14.5
$ws.Range("BD2").Value = 21
$ws.Range("BE2").Value = 60
$ws.Range("BH2").Value = "2026-02-23 04:12:35"
$ws.Range("F2").Value = 2.2
$ws.Range("G2").Value = 2.24
$ws.Range("H2").Value = 3.35
$ws.Range("I2").Value = 3.4
$ws.Range("J2").Value = 3.95
$ws.Range("P2").Value = 2.24
$ws.Range("R2").Value = 1.48
$ws.Range("S2").Value = 2.94
$ws.Range("Z2").Value = 26
$ws.Range("BH3").Value = "2026-02-23 04:12:35"
$ws.Range("BH4").Value = "2026-02-23 04:12:35"
$ws.Range("F4").Value = 15
$ws.Range("K4").Value = 11
$ws.Range("Q4").Value = 1.31
$ws.Range("BH5").Value = "2026-02-23 04:12:35"
$ws.Range("Q5").Value = 2.06
$ws.Range("BH6").Value = "2026-02-23 04:12:35"
$ws.Range("H6").Value = 3.85
$ws.Range("Q6").Value = 2.8
$ws.Range("BH7").Value = "2026-02-23 04:12:35"
$ws.Range("G7").Value = 1.71
$ws.Range("H7").Value = 4.8
$ws.Range("I7").Value = 10.5
$ws.Range("J7").Value = 3.5
$ws.Range("K7").Value = 3.9
$ws.Range("AC8").Value = 970
$ws.Range("AP8").Value = 3.9
$ws.Range("AQ8").Value = 4
$ws.Range("AR8").Value = 5.1
$ws.Range("AS8").Value = 5.8
$ws.Range("AT8").Value = 3.8
$ws.Range("AU8").Value = 3.4
$ws.Range("AV8").Value = 4.6
$ws.Range("AW8").Value = 5.6
$ws.Range("AX8").Value = 4.3
$ws.Range("AY8").Value = 4.1
$ws.Range("AZ8").Value = 4.9
$ws.Range("BA8").Value = 5.7
$ws.Range("BB8").Value = 5.4
$ws.Range("BC8").Value = 5.3
$ws.Range("BD8").Value = 5.6
$ws.Range("BE8").Value = 3.8
$ws.Range("BF8").Value = 5.2
$ws.Range("BG8").Value = 5.7
$ws.Range("BH8").Value = "2026-02-23 04:12:35"
$ws.Range("P8").Value = 1.61
$ws.Range("Q8").Value = 2.32
$ws.Range("BH9").Value = "2026-02-23 04:12:35"
$ws.Range("G9").Value = 3.2
$ws.Range("BH10").Value = "2026-02-23 04:12:35"
$ws.Range("BH11").Value = "2026-02-23 04:12:35"
$ws.Range("H11").Value = 1.98
$ws.Range("I11").Value = 2.18
$ws.Range("BH12").Value = "2026-02-23 04:12:35"
$ws.Range("G12").Value = 2.6
$ws.Range("BH13").Value = "2026-02-23 04:12:35"
$ws.Range("H13").Value = 2.8
$ws.Range("BH14").Value = "2026-02-23 04:12:35"
$ws.Range("BH15").Value = "2026-02-23 04:12:35"
$ws.Range("F15").Value = 1.45
$ws.Range("H15").Value = 1.45
$ws.Range("K15").Value = 3.25
$ws.Range("AP16").Value = 28
$ws.Range("AU16").Value = 11
$ws.Range("BD16").Value = 23
$ws.Range("BH16").Value = "2026-02-23 04:12:35"
$ws.Range("F16").Value = 1.56
$ws.Range("G16").Value = 1.57
$ws.Range("J16").Value = 5.1
$ws.Range("P16").Value = 2.9
$ws.Range("Q16").Value = 1.5
$ws.Range("S16").Value = 2.26
$ws.Range("AT17").Value = 11.5
$ws.Range("BA17").Value = 32
$ws.Range("BD17").Value = 23
$ws.Range("BH17").Value = "2026-02-23 04:12:35"
$ws.Range("G17").Value = 1.52
$ws.Range("P17").Value = 2.84
$ws.Range("T17").Value = 1.71
$ws.Range("AB18").Value = 16.5
$ws.Range("AH18").Value = 26
$ws.Range("AP18").Value = 32
$ws.Range("AR18").Value = 40
$ws.Range("AS18").Value = 18.5
$ws.Range("AU18").Value = 16
$ws.Range("AV18").Value = 30
$ws.Range("AW18").Value = 42
$ws.Range("BH18").Value = "2026-02-23 04:12:35"
$ws.Range("I18").Value = 11.5
$ws.Range("J18").Value = 7
$ws.Range("O18").Value = 1.11
$ws.Range("Q18").Value = 1.34
$ws.Range("R18").Value = 2.1
$ws.Range("S18").Value = 1.84
$ws.Range("T18").Value = 1.69
$ws.Range("U18").Value = 2.38
$ws.Range("BH19").Value = "2026-02-23 04:12:35"
$ws.Range("BH20").Value = "2026-02-23 04:12:35"
$ws.Range("BH21").Value = "2026-02-23 04:12:35"
$ws.Range("F21").Value = 1.89
$ws.Range("BH22").Value = "2026-02-23 04:12:35"
$ws.Range("P22").Value = 1.42
$ws.Range("Q22").Value = 2.92
$ws.Range("BH23").Value = "2026-02-23 04:12:35"
